# Daily data refresh for the "Pais" COVID-19 stats sheet:
#  - bump the "last updated" timestamp
#  - refresh per-country case/death counters
#  - a handful of countries overtook their neighbour in total cases, so the
#    rows for those pairs swap places (country name + stats together)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('A1').Value = 'Datos actualizados a 23 de Septiembre de 2020 a las 15:59'

$ws.Range('B4').Value = 7098766
$ws.Range('C4').Value = 829
$ws.Range('E4').Value = 2545769

$ws.Range('B5').Value = 5664527
$ws.Range('C5').Value = 24031
$ws.Range('D5').Value = 4603424
$ws.Range('E5').Value = 970893
$ws.Range('G5').Value = 189
$ws.Range('H5').Value = 90210

$ws.Range('B19').Value = 331359
$ws.Range('C19').Value = 561
$ws.Range('D19').Value = 313786
$ws.Range('E19').Value = 13004
$ws.Range('G19').Value = 27
$ws.Range('H19').Value = 4569

$ws.Range('B25').Value = 277866
$ws.Range('C25').Value = 690
$ws.Range('E25').Value = 20468
$ws.Range('G25').Value = 7
$ws.Range('H25').Value = 9498

$ws.Range('A51').Value = 'Portugal'
$ws.Range('B51').Value = 70465
$ws.Range('C51').Value = 802
$ws.Range('D51').Value = 46290
$ws.Range('E51').Value = 22247
$ws.Range('G51').Value = 3
$ws.Range('H51').Value = 1928

$ws.Range('A52').Value = 'Etiopia'
$ws.Range('B52').Value = 70422
$ws.Range('D52').Value = 28991
$ws.Range('E52').Value = 40304
$ws.Range('H52').Value = 1127

$ws.Range('E56').Value = 6804
$ws.Range('G56').Value = 4
$ws.Range('H56').Value = 231

$ws.Range('A59').Value = 'Uzbekistan'
$ws.Range('B59').Value = 53275
$ws.Range('C59').Value = 590
$ws.Range('D59').Value = 49649
$ws.Range('E59').Value = 3182
$ws.Range('G59').Value = 2
$ws.Range('H59').Value = 444

$ws.Range('A60').Value = 'Chequia'
$ws.Range('B60').Value = 53158
$ws.Range('C60').Value = 0
$ws.Range('D60').Value = 26151
$ws.Range('E60').Value = 26476
$ws.Range('G60').Value = 0
$ws.Range('H60').Value = 531

$ws.Range('B68').Value = 39524
$ws.Range('C68').Value = 146
$ws.Range('D68').Value = 37122
$ws.Range('E68').Value = 1822
$ws.Range('G68').Value = 2
$ws.Range('H68').Value = 580

$ws.Range('B74').Value = 33080
$ws.Range('C74').Value = 81
$ws.Range('E74').Value = 800
$ws.Range('G74').Value = 1
$ws.Range('H74').Value = 744

$ws.Range('B76').Value = 30097
$ws.Range('C76').Value = 651
$ws.Range('D76').Value = 16430
$ws.Range('E76').Value = 13198
$ws.Range('G76').Value = 9
$ws.Range('H76').Value = 469

$ws.Range('B79').Value = 26081
$ws.Range('C79').Value = 344
$ws.Range('D79').Value = 18634
$ws.Range('E79').Value = 6657
$ws.Range('G79').Value = 12
$ws.Range('H79').Value = 790

$ws.Range('B102').Value = 9475
$ws.Range('C102').Value = 43
$ws.Range('D102').Value = 8244
$ws.Range('E102').Value = 1157
$ws.Range('G102').Value = 1
$ws.Range('H102').Value = 74

$ws.Range('B106').Value = 8646
$ws.Range('C106').Value = 13
$ws.Range('D106').Value = 6551
$ws.Range('E106').Value = 1870
$ws.Range('G106').Value = 2
$ws.Range('H106').Value = 225

$ws.Range('A111').Value = 'Birmania'
$ws.Range('B111').Value = 7177
$ws.Range('C111').Value = 434
$ws.Range('D111').Value = 1951
$ws.Range('E111').Value = 5097
$ws.Range('G111').Value = 14
$ws.Range('H111').Value = 129

$ws.Range('A112').Value = 'Mozambique'
$ws.Range('B112').Value = 7114
$ws.Range('C112').Value = 0
$ws.Range('D112').Value = 4064
$ws.Range('E112').Value = 3005
$ws.Range('G112').Value = 0
$ws.Range('H112').Value = 45

$ws.Range('B120').Value = 5270
$ws.Range('C120').Value = 48
$ws.Range('D120').Value = 4582
$ws.Range('E120').Value = 570
$ws.Range('G120').Value = 1
$ws.Range('H120').Value = 118

$ws.Range('B164').Value = 1337
$ws.Range('C164').Value = 1
$ws.Range('D164').Value = 1219

$ws.Range('A214').Value = 'Islas Malvinas'
$ws.Range('D214').Value = 13
$ws.Range('H214').Value = 0

$ws.Range('A215').Value = 'Montserrat'
$ws.Range('D215').Value = 12
$ws.Range('H215').Value = 1
